$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before U (shifts U:AO -> V:AP), matching the new
# "Enhancement vorhanden" checkbox column added to the template.
$ws.Range("U1").EntireColumn.Insert()

# New header cell for the inserted column.
$ws.Range("U1").Value = "Enhancement vorhanden"

# Re-order the Lateralität / Herdläsion / Zweitläsion headers.
$ws.Range("L1").Value = "Lateralität"
$ws.Range("M1").Value = "Herdläsion"
$ws.Range("N1").Value = "Zweitläsion"

# Approximate the real-Excel "best fit" width the new column ends up with.
$ws.Range("U1").ColumnWidth = 21.25

# Update the UI selection/scroll position to column O (matches the saved view).
$ws.Columns("O").EntireColumn.Select()
